$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

Set-TextValue 'D2' '63.847.48'
$ws.Range('E2').Value = '  +3.25%  '

Set-TextValue 'D3' '3.494.96'
$ws.Range('E3').Value = '  +2.19%  '

$ws.Range('E4').Value = '  -0.03%  '

Set-TextValue 'D5' '414.99'
$ws.Range('E5').Value = '  +1.22%  '

Set-TextValue 'D6' '130.03'
$ws.Range('E6').Value = '  +0.62%  '

Set-TextValue 'D7' '0.633'
$ws.Range('E7').Value = '  +0.31%  '

$ws.Range('E8').Value = '  +0.02%  '

Set-TextValue 'D9' '0.748'
$ws.Range('E9').Value = '  +2.13%  '

Set-TextValue 'D10' '0.158'
$ws.Range('E10').Value = '  +13.51%  '

Set-TextValue 'D11' '42.75'
$ws.Range('E11').Value = '  -1.49%  '

Set-TextValue 'D12' '9.76'
$ws.Range('E12').Value = '  +5.04%  '

Set-TextValue 'D13' '0.0000225'
$ws.Range('E13').Value = '  +0.75%  '

Set-TextValue 'D14' '4.045.20'
$ws.Range('E14').Value = '  +2.13%  '

Set-TextValue 'D15' '0.140'
$ws.Range('E15').Value = '  -0.39%  '

Set-TextValue 'D16' '20.40'
$ws.Range('E16').Value = '  -3.71%  '

Set-TextValue 'D17' '3.500.46'
$ws.Range('E17').Value = '  +2.05%  '

$ws.Range('B18').Value = 'Polygon'
$ws.Range('C18').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D18' '1.10'
$ws.Range('E18').Value = '  +1.58%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D19' '12.49'
$ws.Range('E19').Value = '  +1.15%  '

Set-TextValue 'D20' '63.699.12'
$ws.Range('E20').Value = '  +3.08%  '

Set-TextValue 'D21' '459.89'
$ws.Range('E21').Value = '  -5.72%  '

Set-TextValue 'D22' '90.18'
$ws.Range('E22').Value = '  -1.59%  '

Set-TextValue 'D23' '3.26'
$ws.Range('E23').Value = '  -1.85%  '

Set-TextValue 'D24' '13.34'
$ws.Range('E24').Value = '  -1.25%  '

Set-TextValue 'D25' '10.28'
$ws.Range('E25').Value = '  +10.71%  '

$ws.Range('E26').Value = '  -0.13%  '

Set-TextValue 'D27' '33.68'
$ws.Range('E27').Value = '  -2.19%  '

Set-TextValue 'D28' '4.77'
$ws.Range('E28').Value = '  -0.38%  '

$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D29' '12.70'
$ws.Range('E29').Value = '  +4.73%  '

$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D30' '7.53'
$ws.Range('E30').Value = '  -0.98%  '

$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D31' '2.66'
$ws.Range('E31').Value = '  -1.13%  '

$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D32' '0.170'
$ws.Range('E32').Value = '  +1.45%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D33' '0.114'
$ws.Range('E33').Value = '  -0.62%  '

Set-TextValue 'D34' '40.20'
$ws.Range('E34').Value = '  -4.26%  '

Set-TextValue 'D35' '1.00'
$ws.Range('E35').Value = '  +0.03%  '

Set-TextValue 'D36' '57.59'
$ws.Range('E36').Value = '  -2.06%  '

$ws.Range('E37').Value = '  -1.78%  '

$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D38' '0.0₃0673'
$ws.Range('E38').Value = '  +59.82%  '

$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D39' '3.08'
$ws.Range('E39').Value = '  +5.25%  '

$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D40' '0.999'
$ws.Range('E40').Value = '  +0.11%  '

Set-TextValue 'D41' '2.81'
$ws.Range('E41').Value = '  +2.49%  '

$ws.Range('B42').Value = 'NEARProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D42' '4.61'
$ws.Range('E42').Value = '  +5.74%  '

$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D43' '0.136'
$ws.Range('E43').Value = '  +0.03%  '

$ws.Range('B44').Value = 'LidoDAOToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D44' '3.34'
$ws.Range('E44').Value = '  -3.47%  '

$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D45' '145.67'
$ws.Range('E45').Value = '  -0.41%  '

$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D46' '0.315'
$ws.Range('E46').Value = '  -1.07%  '

Set-TextValue 'D47' '2.01'
$ws.Range('E47').Value = '  -3.68%  '

Set-TextValue 'D48' '2.33'
$ws.Range('E48').Value = '  -1.18%  '

Set-TextValue 'D49' '16.28'
$ws.Range('E49').Value = '  -2.73%  '

Set-TextValue 'D50' '21.79'
$ws.Range('E50').Value = '  -5.46%  '

$ws.Range('E51').Value = '  -1.61%  '
